$wb = $excel.ActiveWorkbook

# Sheet1 = "TC14_CheckOut_PickUp_InStore" (the test-steps sheet), Sheet2 = "Testdata"
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# A new test step ("WAIT") was inserted right before the existing row 13
# ("PRESS_ENTER" on Searchbox), pushing every row from 13..39 down to 14..40.
# Excel's default row-insert duplicates the formatting of the row above, and
# the new Keyword cell (column B) is filled in with the already-existing
# "WAIT" shared string (same one used in row 10).
$ws1.Rows("13:13").Insert()
$ws1.Range("B13").Value = "WAIT"

# Cosmetic: restore the selections left behind in the saved file. Sheet2's
# selection/view is updated first so that Sheet1 ends up as the active
# (tabSelected) sheet, matching the original file.
[void]$ws2.Activate()
[void]$ws2.Range("B11").Select()

[void]$ws1.Activate()
[void]$ws1.Range("D12").Select()
